$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 23.47930048478094
$ws.Range("C2").Value = 6.270156305745203
$ws.Range("D2").Value = 6.103539760974283
$ws.Range("E2").Value = 8.903859682885081
$ws.Range("G2").Value = 3.811569418557166
$ws.Range("K2").Value = 19.78711842111927
$ws.Range("L2").Value = 10.26232679005846
$ws.Range("M2").Value = 19.20166858335622
$ws.Range("N2").Value = 29.08391994935222
$ws.Range("B3").Value = 23.33036574243604
$ws.Range("C3").Value = 6.061030079424166
$ws.Range("D3").Value = 6.004260750733199
$ws.Range("E3").Value = 8.91558337082296
$ws.Range("G3").Value = 3.815654318945504
$ws.Range("K3").Value = 19.67727322407919
$ws.Range("L3").Value = 10.2765757029304
$ws.Range("M3").Value = 19.19455112777568
$ws.Range("N3").Value = 28.9987185483454
$ws.Range("B4").Value = 23.24464556243133
$ws.Range("C4").Value = 5.931115120121863
$ws.Range("D4").Value = 5.944429068301987
$ws.Range("E4").Value = 8.923304402984497
$ws.Range("G4").Value = 3.818290822549458
$ws.Range("K4").Value = 19.61463792043736
$ws.Range("L4").Value = 10.28678776082565
$ws.Range("M4").Value = 19.1939014786496
$ws.Range("N4").Value = 28.94726111149516
$ws.Range("B5").Value = 23.21118167481673
$ws.Range("C5").Value = 5.877887677120825
$ws.Range("D5").Value = 5.920363248712982
$ws.Range("E5").Value = 8.926582510885419
$ws.Range("G5").Value = 3.81939762398707
$ws.Range("K5").Value = 19.59034183236233
$ws.Range("L5").Value = 10.29131715861949
$ws.Range("M5").Value = 19.19457301060381
$ws.Range("N5").Value = 28.92651148922774
$ws.Range("B6").Value = 23.20571443571866
$ws.Range("C6").Value = 5.869034893666631
$ws.Range("D6").Value = 5.916387179585785
$ws.Range("E6").Value = 8.927134802984869
$ws.Range("G6").Value = 3.819583368452435
$ws.Range("K6").Value = 19.58638218203081
$ws.Range("L6").Value = 10.2920914804902
$ws.Range("M6").Value = 19.19474107813798
$ws.Range("N6").Value = 28.92307942092209
$ws.Range("B7").Value = 23.24418827893558
$ws.Range("C7").Value = 5.93039830558448
$ws.Range("D7").Value = 5.944103185314913
$ws.Range("E7").Value = 8.923348078929358
$ws.Range("G7").Value = 3.818305617886984
$ws.Range("K7").Value = 19.61430525825638
$ws.Range("L7").Value = 10.286847356319
$ws.Range("M7").Value = 19.1939067436922
$ws.Range("N7").Value = 28.94698037988569
$ws.Range("B8").Value = 23.42677674692731
$ws.Range("C8").Value = 6.198417331548114
$ws.Range("D8").Value = 6.069094545139802
$ws.Range("E8").Value = 8.907793717721338
$ws.Range("G8").Value = 3.81295132873857
$ws.Range("K8").Value = 19.74825731022045
$ws.Range("L8").Value = 10.26693613069139
$ws.Range("M8").Value = 19.1984428596143
$ws.Range("N8").Value = 29.05436736383136
$ws.Range("B9").Value = 23.82892900692275
$ws.Range("C9").Value = 6.70824116032681
$ws.Range("D9").Value = 6.321673482755466
$ws.Range("E9").Value = 8.881425187083217
$ws.Range("G9").Value = 3.803464203120158
$ws.Range("K9").Value = 20.04818919776145
$ws.Range("L9").Value = 10.23950034296067
$ws.Range("M9").Value = 19.23681007067914
$ws.Range("N9").Value = 29.27155722605797
$ws.Range("B10").Value = 24.14932651110167
$ws.Range("C10").Value = 7.068651919870724
$ws.Range("D10").Value = 6.509840291946685
$ws.Range("E10").Value = 8.864553976506517
$ws.Range("G10").Value = 3.79710314476493
$ws.Range("K10").Value = 20.28991607105822
$ws.Range("L10").Value = 10.22642082134081
$ws.Range("M10").Value = 19.28287525371012
$ws.Range("N10").Value = 29.43489156970688
$ws.Range("B11").Value = 24.30004302761672
$ws.Range("C11").Value = 7.22872128484396
$ws.Range("D11").Value = 6.595599579088408
$ws.Range("E11").Value = 8.857418262385943
$ws.Range("G11").Value = 3.79433985044447
$ws.Range("K11").Value = 20.40420415223384
$ws.Range("L11").Value = 10.22200654489918
$ws.Range("M11").Value = 19.30768414060736
$ws.Range("N11").Value = 29.50995446039576
$ws.Range("B12").Value = 24.35778882418916
$ws.Range("C12").Value = 7.288718381245256
$ws.Range("D12").Value = 6.628064781246357
$ws.Range("E12").Value = 8.854793380301277
$ws.Range("G12").Value = 3.793312076743179
$ws.Range("K12").Value = 20.4480750519924
$ws.Range("L12").Value = 10.2205556197104
$ws.Range("M12").Value = 19.31762938399219
$ws.Range("N12").Value = 29.53848265708938
$ws.Range("B13").Value = 24.34532298660494
$ws.Range("C13").Value = 7.275825278987223
$ws.Range("D13").Value = 6.621073778667632
$ws.Range("E13").Value = 8.855355263683846
$ws.Range("G13").Value = 3.793532599762405
$ws.Range("K13").Value = 20.43860079953537
$ws.Range("L13").Value = 10.22085829136969
$ws.Range("M13").Value = 19.31546306812873
$ws.Range("N13").Value = 29.53233409291134
$ws.Range("B14").Value = 24.30478056828589
$ws.Range("C14").Value = 7.233669977677263
$ws.Range("D14").Value = 6.598270896669469
$ws.Range("E14").Value = 8.857200764876563
$ws.Range("G14").Value = 3.794254922296801
$ws.Range("K14").Value = 20.40780173628148
$ws.Range("L14").Value = 10.22188275505269
$ws.Range("M14").Value = 19.30849132754583
$ws.Range("N14").Value = 29.51229945269531
$ws.Range("B15").Value = 24.28003352447797
$ws.Range("C15").Value = 7.207766549794684
$ws.Range("D15").Value = 6.584301208715261
$ws.Range("E15").Value = 8.858341241216426
$ws.Range("G15").Value = 3.79469978769967
$ws.Range("K15").Value = 20.38901268473925
$ws.Range("L15").Value = 10.22253900013481
$ws.Range("M15").Value = 19.30429253480645
$ws.Range("N15").Value = 29.5000409566431
$ws.Range("B16").Value = 24.13957325412848
$ws.Range("C16").Value = 7.058107647991796
$ws.Range("D16").Value = 6.504236104985186
$ws.Range("E16").Value = 8.865031137504412
$ws.Range("G16").Value = 3.797286345661573
$ws.Range("K16").Value = 20.28253161840624
$ws.Range("L16").Value = 10.22674019357983
$ws.Range("M16").Value = 19.28133119231799
$ws.Range("N16").Value = 29.43000084710001
$ws.Range("B17").Value = 24.05464833623587
$ws.Range("C17").Value = 6.96525783647204
$ws.Range("D17").Value = 6.455136112834321
$ws.Range("E17").Value = 8.869273057202456
$ws.Range("G17").Value = 3.798906421752133
$ws.Range("K17").Value = 20.21829614680059
$ws.Range("L17").Value = 10.22971071175424
$ws.Range("M17").Value = 19.2682298310124
$ws.Range("N17").Value = 29.38722431271903
$ws.Range("B18").Value = 24.00627177099962
$ws.Range("C18").Value = 6.911490672647642
$ws.Range("D18").Value = 6.426912194981321
$ws.Range("E18").Value = 8.871763653789507
$ws.Range("G18").Value = 3.79985052627328
$ws.Range("K18").Value = 20.18175856743361
$ws.Range("L18").Value = 10.23156383209623
$ws.Range("M18").Value = 19.2610572347641
$ws.Range("N18").Value = 29.36269240657137
$ws.Range("B19").Value = 23.98997425982137
$ws.Range("C19").Value = 6.893225707260925
$ws.Range("D19").Value = 6.417360026912472
$ws.Range("E19").Value = 8.872615652524967
$ws.Range("G19").Value = 3.800172296628933
$ws.Range("K19").Value = 20.16945869006825
$ws.Range("L19").Value = 10.23221609862784
$ws.Range("M19").Value = 19.2586911536823
$ws.Range("N19").Value = 29.35439887332591
$ws.Range("B20").Value = 24.06364039109334
$ws.Range("C20").Value = 6.975179815908662
$ws.Range("D20").Value = 6.460361364990454
$ws.Range("E20").Value = 8.868816246405185
$ws.Range("G20").Value = 3.798732691887275
$ws.Range("K20").Value = 20.22509200602188
$ws.Range("L20").Value = 10.22937953494384
$ws.Range("M20").Value = 19.26958695196069
$ws.Range("N20").Value = 29.39177053606426
$ws.Range("B21").Value = 24.31667093796241
$ws.Range("C21").Value = 7.246069215293828
$ws.Range("D21").Value = 6.604969185288823
$ws.Range("E21").Value = 8.856656601748298
$ws.Range("G21").Value = 3.794042254136428
$ws.Range("K21").Value = 20.41683232858318
$ws.Range("L21").Value = 10.22157585797832
$ws.Range("M21").Value = 19.31052418236493
$ws.Range("N21").Value = 29.51818135146034
$ws.Range("B22").Value = 24.48594213360071
$ws.Range("C22").Value = 7.419487775573384
$ws.Range("D22").Value = 6.699406059700008
$ws.Range("E22").Value = 8.849159771013889
$ws.Range("G22").Value = 3.791085290994246
$ws.Range("K22").Value = 20.54558471645431
$ws.Range("L22").Value = 10.21776179982134
$ws.Range("M22").Value = 19.34048670660155
$ws.Range("N22").Value = 29.60139855322597
$ws.Range("B23").Value = 24.39525594852129
$ws.Range("C23").Value = 7.327280451512906
$ws.Range("D23").Value = 6.64902062509962
$ws.Range("E23").Value = 8.853119862799726
$ws.Range("G23").Value = 3.792653589842034
$ws.Range("K23").Value = 20.47656250377861
$ws.Range("L23").Value = 10.21967982321828
$ws.Range("M23").Value = 19.3242029215764
$ws.Range("N23").Value = 29.55693104788145
$ws.Range("B24").Value = 24.05957368485607
$ws.Range("C24").Value = 6.970695288515977
$ws.Range("D24").Value = 6.45799901250062
$ws.Range("E24").Value = 8.869022609013257
$ws.Range("G24").Value = 3.798811195588427
$ws.Range("K24").Value = 20.2220183740775
$ws.Range("L24").Value = 10.22952880729802
$ws.Range("M24").Value = 19.26897227740075
$ws.Range("N24").Value = 29.38971499837118
$ws.Range("B25").Value = 23.71560802608567
$ws.Range("C25").Value = 6.572505401402981
$ws.Range("D25").Value = 6.252747682948763
$ws.Range("E25").Value = 8.888117941747216
$ws.Range("G25").Value = 3.805923161909333
$ws.Range("K25").Value = 19.96319357976238
$ws.Range("L25").Value = 10.24567935273459
$ws.Range("M25").Value = 19.22328288217824
$ws.Range("N25").Value = 29.21212062681202
